# Refresh cryptocurrency price/volume snapshot (cryptos.xlsx) with latest values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.645.00"
$ws.Range("E2").Value = "  -0.73%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.530.28"
$ws.Range("E3").Value = "  -2.74%  "
$ws.Range("E4").Value = "  +0.18%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "201.62"
$ws.Range("E5").Value = "  +2.56%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "552.95"
$ws.Range("E6").Value = "  -3.73%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.519.51"
$ws.Range("E7").Value = "  -2.87%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.612"
$ws.Range("E8").Value = "  -1.15%  "
$ws.Range("E9").Value = "  +0.01%  "
$ws.Range("E10").Value = "  -2.93%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "62.64"
$ws.Range("E11").Value = "  +11.29%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.144"
$ws.Range("E12").Value = "  -7.08%  "
$ws.Range("E13").Value = "  -7.91%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.88"
$ws.Range("E14").Value = "  -1.89%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.088.68"
$ws.Range("E15").Value = "  -2.97%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.528.18"
$ws.Range("E16").Value = "  -2.84%  "
$ws.Range("E17").Value = "  -1.73%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "18.50"
$ws.Range("E18").Value = "  -0.23%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "67.490.70"
$ws.Range("E19").Value = "  -0.95%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.86"
$ws.Range("E20").Value = "  -5.26%  "
$ws.Range("E21").Value = "  -5.46%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "394.50"
$ws.Range("E22").Value = "  -2.09%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.08"
$ws.Range("E23").Value = "  -9.13%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.00"
$ws.Range("E24").Value = "  -5.69%  "
$ws.Range("E25").Value = "  -2.98%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.89"
$ws.Range("E26").Value = "  +0.34%  "
$ws.Range("B27").Value = "InternetComputer(DFINITY)"
$ws.Range("C27").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.28"
$ws.Range("E27").Value = "  -2.87%  "
$ws.Range("B28").Value = "ImmutableX"
$ws.Range("C28").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.83"
$ws.Range("E28").Value = "  -4.42%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.88"
$ws.Range("E29").Value = "  -3.45%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "717.17"
$ws.Range("E30").Value = "  +3.89%  "
$ws.Range("E31").Value = "  -2.08%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.12"
$ws.Range("E32").Value = "  -13.62%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "11.78"
$ws.Range("E33").Value = "  -3.59%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "63.67"
$ws.Range("E34").Value = "  -1.70%  "
$ws.Range("E35").Value = "  -4.89%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "38.63"
$ws.Range("E36").Value = "  -9.49%  "
$ws.Range("E37").Value = "  -0.04%  "
$ws.Range("E38").Value = "  -6.63%  "
$ws.Range("B39").Value = "Kaspa"
$ws.Range("C39").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.133"
$ws.Range("E39").Value = "  -4.44%  "
$ws.Range("B40").Value = "ThetaToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.05"
$ws.Range("E40").Value = "  -3.85%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.080.72"
$ws.Range("E41").Value = "  -4.67%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.998"
$ws.Range("E42").Value = "  -0.17%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0₃0683"
$ws.Range("E43").Value = "  -13.32%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.55"
$ws.Range("E44").Value = "  -12.72%  "
$ws.Range("B45").Value = "WEMIXToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.78"
$ws.Range("E45").Value = "  +6.13%  "
$ws.Range("B46").Value = "VeChain"
$ws.Range("C46").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0412"
$ws.Range("E46").Value = "  -2.06%  "
$ws.Range("B47").Value = "dogwifhat"
$ws.Range("C47").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.71"
$ws.Range("E47").Value = "  -13.11%  "
$ws.Range("E48").Value = "  -3.06%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "137.86"
$ws.Range("E49").Value = "  -3.34%  "
$ws.Range("E50").Value = "  -7.18%  "
$ws.Range("E51").Value = "  -6.84%  "
